$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H74").Value = 3899
$ws.Range("I74").Value = 3770.4285
$ws.Range("J74").Value = 4199
$ws.Range("K74").Value = 3770.4285
$ws.Range("L74").Value = 4199
$ws.Range("M74").Value = -2834.4285
$ws.Range("N74").Value = -6071

$ws.Range("H76").Value = 1956193
$ws.Range("I76").Value = 2605413.5
$ws.Range("J76").Value = 8531.666999999999
$ws.Range("K76").Value = 2605413.5
$ws.Range("L76").Value = 8531.666999999999
$ws.Range("M76").Value = -2605098.5
$ws.Range("N76").Value = -9161.666999999999

$ws.Range("H77").Value = 3899
$ws.Range("I77").Value = 3770.4285
$ws.Range("J77").Value = 4199
$ws.Range("K77").Value = 18852.1425
$ws.Range("L77").Value = 20995
$ws.Range("M77").Value = -14172.1425
$ws.Range("N77").Value = -30355

$ws.Range("H79").Value = 1956193
$ws.Range("I79").Value = 2605413.5
$ws.Range("J79").Value = 8531.666999999999
$ws.Range("K79").Value = 2605413.5
$ws.Range("L79").Value = 8531.666999999999
$ws.Range("M79").Value = -2604321.5
$ws.Range("N79").Value = -10715.667

$ws.Range("H120").Value = 33165
$ws.Range("J120").Value = 33165
$ws.Range("L120").Value = 33165
$ws.Range("N120").Value = -42841

$ws.Range("H137").Value = 1922
$ws.Range("I137").Value = 1236.125
$ws.Range("J137").Value = 2344.077
$ws.Range("K137").Value = 3708.375
$ws.Range("L137").Value = 7032.231000000001
$ws.Range("M137").Value = -1158.375
$ws.Range("N137").Value = -12132.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 399.5
$ws.Range("I5").Value = 399.5
$ws.Range("K5").Value = 399.5
$ws.Range("M5").Value = -287.5

$ws.Range("H45").Value = 1893.8
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 399.5
$ws.Range("I4").Value = 399.5
$ws.Range("K4").Value = 399.5
$ws.Range("M4").Value = -284.5

$ws.Range("H86").Value = 144809.36
$ws.Range("I86").Value = 1864
$ws.Range("K86").Value = 1864
$ws.Range("M86").Value = -741

$ws.Range("H89").Value = 144809.36
$ws.Range("I89").Value = 1864
$ws.Range("K89").Value = 9320
$ws.Range("M89").Value = -3704

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 15181.637
$ws.Range("I97").Value = 7249.5
$ws.Range("K97").Value = 7249.5
$ws.Range("M97").Value = -6258.5

$ws.Range("H105").Value = 1787.5714
$ws.Range("I105").Value = 1899.9412
$ws.Range("K105").Value = 1899.9412
$ws.Range("M105").Value = -152.9412

$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

$ws.Range("H134").Value = 10211.759
$ws.Range("I134").Value = 10796.333
$ws.Range("J134").Value = 7405.8
$ws.Range("K134").Value = 32388.999
$ws.Range("L134").Value = 22217.4
$ws.Range("M134").Value = -29853.999
$ws.Range("N134").Value = -27287.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.25
$ws.Range("I7").Value = 70.666664
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 70.666664
$ws.Range("L7").Value = 65
$ws.Range("M7").Value = 42.333336
$ws.Range("N7").Value = -291

$ws.Range("H16").Value = 889.5
$ws.Range("I16").Value = 790
$ws.Range("K16").Value = 790
$ws.Range("M16").Value = -503

$ws.Range("H22").Value = 667.36365
$ws.Range("I22").Value = 259
$ws.Range("J22").Value = 1007.6667
$ws.Range("K22").Value = 259
$ws.Range("L22").Value = 1007.6667
$ws.Range("M22").Value = 91
$ws.Range("N22").Value = -1707.6667

$ws.Range("H31").Value = 2541.6924
$ws.Range("I31").Value = 1503.8182
$ws.Range("J31").Value = 8250
$ws.Range("K31").Value = 1503.8182
$ws.Range("L31").Value = 8250
$ws.Range("M31").Value = -1208.8182
$ws.Range("N31").Value = -8840

$ws.Range("H34").Value = 2541.6924
$ws.Range("I34").Value = 1503.8182
$ws.Range("J34").Value = 8250
$ws.Range("K34").Value = 1503.8182
$ws.Range("L34").Value = 8250
$ws.Range("M34").Value = -1301.8182
$ws.Range("N34").Value = -8654

$ws.Range("H113").Value = 889.5
$ws.Range("I113").Value = 790
$ws.Range("K113").Value = 790
$ws.Range("M113").Value = 1380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 777.59
$ws.Range("J131").Value = 814.2717
$ws.Range("L131").Value = 2442.8151
$ws.Range("N131").Value = -12522.8151

$ws.Range("H132").Value = 1466.6666
$ws.Range("I132").Value = 1463.6364
$ws.Range("K132").Value = 13172.7276
$ws.Range("M132").Value = -10642.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H55").Value = 526
$ws.Range("I55").Value = 489.5
$ws.Range("K55").Value = 489.5
$ws.Range("M55").Value = -316.5

$ws.Range("H122").Value = 3791.6924
$ws.Range("J122").Value = 4675
$ws.Range("L122").Value = 14025
$ws.Range("N122").Value = -18925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 806.625
$ws.Range("I107").Value = 441.66666
$ws.Range("K107").Value = 1324.99998
$ws.Range("M107").Value = 595.0000199999999

$ws.Range("H132").Value = 4727.729
$ws.Range("I132").Value = 971.6087
$ws.Range("K132").Value = 2914.8261
$ws.Range("M132").Value = -384.8261000000002
